$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expand the workers table -------------------------------------------
# The table currently has 3 data rows: 16 and 17 use the "middle" row style,
# 18 uses the "last" row style (bottom border closing the table). The new
# statement has 12 data rows: 11 "middle" style rows (16-26) followed by a
# single "last" style row (27). Insert 9 blank rows above the current last
# row (row 18) so that row shifts down to become the new row 27, then copy
# the "middle" row formatting (taken from row 17) onto the newly inserted
# rows 18-26.
$ws.Rows("18:26").Insert() | Out-Null

$ws.Range("B17:J17").Copy() | Out-Null
$ws.Range("B18:J26").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Fill in the new worker / period / value data -----------------------
$data = @(
    @{Row=16; Doc="37557547";   Name="ROSA ELENA PALOMINO LIZARAZO"; Period="2507"; Mora=64000;  Salario=1600000},
    @{Row=17; Doc="37557547";   Name="ROSA ELENA PALOMINO LIZARAZO"; Period="2506"; Mora=64000;  Salario=1600000},
    @{Row=18; Doc="9116913";    Name="ISAAC MANUEL RODELO JIMENEZ";  Period="2506"; Mora=76000;  Salario=1900000},
    @{Row=19; Doc="9116913";    Name="ISAAC MANUEL RODELO JIMENEZ";  Period="2505"; Mora=76000;  Salario=1900000},
    @{Row=20; Doc="9116913";    Name="ISAAC MANUEL RODELO JIMENEZ";  Period="2504"; Mora=76000;  Salario=1900000},
    @{Row=21; Doc="1104374560"; Name="JAIRO MANUEL ZABALETA TIRADO"; Period="2505"; Mora=82000;  Salario=2050000},
    @{Row=22; Doc="9691206";    Name="ROBINSON CASTRO CHOGO";        Period="2506"; Mora=114000; Salario=2850000},
    @{Row=23; Doc="9096340";    Name="ALEXANDER SIERRA ARIAS";       Period="2504"; Mora=70000;  Salario=1750000},
    @{Row=24; Doc="1065866974"; Name="JHON FREIDER MENESES NAVARRO"; Period="2506"; Mora=56940;  Salario=1423500},
    @{Row=25; Doc="1065866974"; Name="JHON FREIDER MENESES NAVARRO"; Period="2505"; Mora=56940;  Salario=1423500},
    @{Row=26; Doc="1065866974"; Name="JHON FREIDER MENESES NAVARRO"; Period="2504"; Mora=49348;  Salario=1423500},
    @{Row=27; Doc="1104378817"; Name="JORGE LUIS ZAPATA LARA";       Period="2504"; Mora=70000;  Salario=1750000}
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Range("B$r").Value2 = "CC"
    $ws.Range("C$r").Value2 = $row.Doc
    $ws.Range("D$r").Value2 = $row.Name
    $ws.Range("E$r").Value2 = $row.Period
    $ws.Range("F$r").Value2 = $row.Mora
    $ws.Range("G$r").Value2 = $row.Salario
}

# --- Update the summary figures above the table --------------------------
$ws.Range("E11").Value2 = 855228   # VALOR MORA total
$ws.Range("C13").Value2 = 7        # Cant. Trabajadores
$ws.Range("F13").Value2 = 4        # Cant. Periodos
